$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new "Wins" / "Losses" / "Ties" columns, styled like the
# existing header cells (bold, centered, top-aligned, thin border).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$headerRange = $ws.Range("AC1:AE1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Season record: Wins / Losses / Ties for every player row.
$wins = 88
$losses = 74
$ties = 0

$lastRow = 40
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins
    $ws.Cells.Item($row, 30).Value = $losses
    $ws.Cells.Item($row, 31).Value = $ties
}
